$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: E2 price input changes from 50 to 62 ---
$ws1.Range("E2").Value = 62
# E2 picks up the 2-decimal numeric format already used by F10:F12
$ws1.Range("F10").Copy()
$ws1.Range("E2").PasteSpecial(-4122)

# New formatted-but-empty cell F2 (percentage style, matching Z29's prior style)
$ws2.Range("Z29").Copy()
$ws1.Range("F2").PasteSpecial(-4122)

$ws1.Range("E2").Select()

# --- Sheet2: assumption changes ---
# ROIC 6% -> 4%
$ws2.Range("Z29").Value = 0.04
# Discount rate 8.5% -> 9%
$ws2.Range("Z31").Value = 0.09

# ROIC & Maturity now share Discount's "0.0%" number format
$ws2.Range("Z31").Copy()
$ws2.Range("Z29:Z30").PasteSpecial(-4122)

# NPV/share (Z33) switches from currency format to plain 2-decimal format
$ws1.Range("F10").Copy()
$ws2.Range("Z33").PasteSpecial(-4122)

$ws2.Range("Z32").Select()
